$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.780.68"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.311.63"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.30"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.35"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.15"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.93"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0782"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").Value = "2.671.74"
$ws.Range("D16").Value = "2.293.68"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "42.727.23"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.20"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.67"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.28"
$ws.Range("E23").Value = "  +6.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.03"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.28"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("E28").Value = "  +15.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.35"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.11"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.10"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.61"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.85"
$ws.Range("E42").Value = "  +14.24%  "
$ws.Range("D43").Value = "1.925.76"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0279"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.11"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.73"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").Value = "2.540.68"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.37"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  +1.94%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
